$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("I").ColumnWidth = 28.5703125
$ws.Columns("J").ColumnWidth = 32.140625
$ws.Columns("K").ColumnWidth = 28
$ws.Columns("L").ColumnWidth = 37.42578125
$ws.Columns("M").ColumnWidth = 34
$ws.Columns("N").ColumnWidth = 21
